$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 120 (the "FAO" row), shifting that
# row and everything below it down by one.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row with the new abbreviation entry.
$ws.Range("A120").Value = "FATF"
$ws.Range("B120").Value = "Financial Action Task Force"
$ws.Range("C120").Value = "Financial Action Task Force"

# Match the formatting of the surrounding data rows (the insert operation
# can leave the new row with a slightly different style than its neighbours).
$ws.Range("A121:C121").Copy()
$ws.Range("A120:C120").PasteSpecial(-4122)
$excel.CutCopyMode = $false
